# chore: adapt column header formatting to respective input file names
#
# 1. Rename the "_old" / "_new" header-name suffixes to "_FV2210" / "_FV2304"
#    respectively (leaving the "diff" column header untouched).
# 2. Turn the data range into a real Excel Table ("Table1") with an
#    autofilter.
# 3. Freeze the header row (row 1) in the sheet view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header cells -------------------------------------------------
$usedRange = $ws.UsedRange
$lastCol = $usedRange.Columns.Count

for ($c = 1; $c -le $lastCol; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $headerText = $cell.Value2
    if ($headerText -ne $null) {
        if ($headerText.EndsWith("_old")) {
            $cell.Value2 = $headerText.Substring(0, $headerText.Length - 4) + "_FV2210"
        } elseif ($headerText.EndsWith("_new")) {
            $cell.Value2 = $headerText.Substring(0, $headerText.Length - 4) + "_FV2304"
        }
    }
}

# --- 2. Convert the used range into a table ---------------------------------
$dataRange = $ws.UsedRange
$lo = $ws.ListObjects.Add(
    [Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange,
    $dataRange,
    $null,
    [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes
)
$lo.Name = "Table1"

# --- 3. Freeze the header row ------------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
